# Generate Report for Handoff
# Adds a new localization-status row for 468b9ed0-1b47-45f9-a922-0e8b8714b07e.md
# to the Overview, zh-cn and de-de sheets (mirrors the existing
# 0a7476a9-66c5-4280-8d56-91b9965d7921.md row on each sheet).

$wb = $excel.ActiveWorkbook

$commit = "9d70a0a06bdb2103140d0fda9bd3e9f0707b0c30"
$newFile = "468b9ed0-1b47-45f9-a922-0e8b8714b07e.md"
$newFileDisplay = "e2e\468b9ed0-1b47-45f9-a922-0e8b8714b07e.md"
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newFile"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = $newFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", $newFileDisplay)
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 4).Value = ""
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = "2016-08-29 08:42:47"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, "", "", $newFile)
$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = "'False"
$wsZhCn.Cells.Item(3, 7).Value = "468b9ed0-1b47-45f9-a922-0e8b8714b07e.4d089905522975fa48e647c5cdcab7729697423d.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 8).Value = "2016-08-29 08:42:42"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 9).Value = ""
$wsZhCn.Cells.Item(3, 10).Value = ""
$wsZhCn.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 12).Value = ""
$wsZhCn.Cells.Item(3, 13).Value = "'True"
$wsZhCn.Cells.Item(3, 14).Value = ""
$wsZhCn.Cells.Item(3, 15).Value = "'False"
$wsZhCn.Cells.Item(3, 16).Value = ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, "", "", $newFile)
$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = "'False"
$wsDeDe.Cells.Item(3, 7).Value = "468b9ed0-1b47-45f9-a922-0e8b8714b07e.4d089905522975fa48e647c5cdcab7729697423d.de-de.xlf"
$wsDeDe.Cells.Item(3, 8).Value = "2016-08-29 08:42:47"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 9).Value = ""
$wsDeDe.Cells.Item(3, 10).Value = ""
$wsDeDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 12).Value = ""
$wsDeDe.Cells.Item(3, 13).Value = "'True"
$wsDeDe.Cells.Item(3, 14).Value = ""
$wsDeDe.Cells.Item(3, 15).Value = "'False"
$wsDeDe.Cells.Item(3, 16).Value = ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
